$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(15)

$urls = @(
    "https://japan.cnet.com/release/30409909/?ref=rss",
    "https://japan.cnet.com/release/30410024/?ref=rss",
    "https://japan.cnet.com/release/30409893/?ref=rss",
    "https://japan.cnet.com/release/30409855/?ref=rss",
    "https://japan.cnet.com/release/30410063/?ref=rss",
    "https://japan.cnet.com/release/30409892/?ref=rss",
    "https://japan.cnet.com/release/30410098/",
    "https://japan.cnet.com/release/30410048/?ref=rss",
    "https://japan.cnet.com/release/30409583/",
    "https://japan.cnet.com/release/30409727/",
    "https://japan.cnet.com/release/30409394/",
    "https://japan.cnet.com/release/30409148/",
    "https://japan.cnet.com/release/30409794/",
    "https://japan.cnet.com/release/30408791/",
    "https://japan.cnet.com/release/30408925/",
    "https://japan.cnet.com/release/30407868/",
    "https://japan.cnet.com/article/35146156/",
    "https://japan.cnet.com/release/30405956/",
    "https://japan.cnet.com/release/30405909/",
    "https://japan.cnet.com/release/30405516/"
)

$row = 3
foreach ($url in $urls) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $url
    $cell.Style = "Normal"
    $row = $row + 1
}
